# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Decision-table rule row 11 (column B, the "Rule" name) is renamed from
# "R40" to "1" on the Rules sheet.  It must remain a text value (it is
# stored as a shared string, not a number) even though it looks numeric,
# so the cell is entered with a leading apostrophe to force text, exactly
# like typing '1 into Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("B11").Value = "'1"
